$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws1.Cells.Item(4, 6).Value = 1283
$ws1.Cells.Item(6, 6).Value = 311
$ws1.Cells.Item(7, 6).Value = 1117
$ws1.Cells.Item(8, 6).Value = 431
$ws1.Cells.Item(9, 6).Value = 6956
$ws1.Cells.Item(12, 6).Value = 2033
$ws1.Cells.Item(13, 6).Value = 7851
$ws1.Cells.Item(15, 6).Value = 49
$ws1.Cells.Item(16, 6).Value = 5452
$ws1.Cells.Item(17, 6).Value = 44
$ws1.Cells.Item(18, 6).Value = 2328
$ws1.Cells.Item(19, 6).Value = 980
$ws1.Cells.Item(21, 6).Value = 275
$ws1.Cells.Item(22, 6).Value = 370
$ws1.Cells.Item(25, 6).Value = 320
$ws1.Cells.Item(26, 6).Value = 236
$ws1.Cells.Item(28, 6).Value = 2109
$ws1.Cells.Item(30, 6).Value = 239
$ws1.Cells.Item(32, 6).Value = 53
$ws1.Cells.Item(33, 6).Value = 548
$ws1.Cells.Item(36, 6).Value = 1421
$ws1.Cells.Item(37, 6).Value = 26
$ws1.Cells.Item(39, 6).Value = 2156

$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws2.Cells.Item(4, 6).Value = 37

$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws3.Cells.Item(3, 6).Value = 1262

$ws4 = $wb.Worksheets.Item(4)  # 全部类型
$ws4.Cells.Item(4, 6).Value = 1262
$ws4.Cells.Item(6, 6).Value = 1283
$ws4.Cells.Item(9, 6).Value = 311
$ws4.Cells.Item(10, 6).Value = 1117
$ws4.Cells.Item(11, 6).Value = 431
$ws4.Cells.Item(12, 6).Value = 6956
$ws4.Cells.Item(15, 6).Value = 2033
$ws4.Cells.Item(16, 6).Value = 7851
$ws4.Cells.Item(18, 6).Value = 49
$ws4.Cells.Item(19, 6).Value = 5452
$ws4.Cells.Item(20, 6).Value = 44
$ws4.Cells.Item(21, 6).Value = 2328
$ws4.Cells.Item(22, 6).Value = 980
$ws4.Cells.Item(24, 6).Value = 275
$ws4.Cells.Item(25, 6).Value = 370
$ws4.Cells.Item(29, 6).Value = 37
$ws4.Cells.Item(30, 6).Value = 320
$ws4.Cells.Item(31, 6).Value = 236
$ws4.Cells.Item(33, 6).Value = 2109
$ws4.Cells.Item(35, 6).Value = 239
$ws4.Cells.Item(37, 6).Value = 53
$ws4.Cells.Item(38, 6).Value = 548
$ws4.Cells.Item(42, 6).Value = 1421
$ws4.Cells.Item(43, 6).Value = 26
$ws4.Cells.Item(45, 6).Value = 2157

$wb.Save()
